$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.95"
$ws.Range("D3").Value = "'23.10"
$ws.Range("D4").Value = "'5.411"
$ws.Range("D5").Value = "'0.05895"
$ws.Range("D7").Value = "'6.542"
$ws.Range("D8").Value = "'0.8103"
$ws.Range("D9").Value = "'0.9354"
$ws.Range("D11").Value = "'0.07414"
$ws.Range("D13").Value = "'0.03040"
$ws.Range("D14").Value = "'0.09354"
$ws.Range("D15").Value = "'3.854"
$ws.Range("D16").Value = "'0.001586"
$ws.Range("D17").Value = "'0.04681"
$ws.Range("D18").Value = "'0.0005908"
$ws.Range("D19").Value = "'0.005870"
$ws.Range("D20").Value = "'0.001266"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("D21").Value = "'0.004904"
$ws.Range("D22").Value = "'0.00006804"
$ws.Range("D24").Value = "'2.110"
$ws.Range("D26").Value = "'0.1330"
$ws.Range("D27").Value = "'0.0002285"
$ws.Range("D40").Value = "'0.03969"
$ws.Range("D41").Value = "'0.006180"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.003001"
$ws.Range("D44").Value = "'0.009720"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").Value = "'0.00005190"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.6703"
$ws.Range("D48").Value = "'0.002388"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002001"
